$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; this shifts existing rows 47..149 down to 48..150
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with its data
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44973
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112022
$ws.Cells.Item(47, 7).Value = "Arveja Verde"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 15
$ws.Cells.Item(47, 11).Value = 35000
$ws.Cells.Item(47, 12).Value = 35000
$ws.Cells.Item(47, 13).Value = 35000
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 16).Value = 1400
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
